$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row updates
$ws.Range("A1").Value = "RTDO L"
$ws.Range("D1").Value = "RTDO V"
$ws.Range("E1").Value = "Jornada"
$ws.Range("F1").Value = "RTDO L.1"
$ws.Range("I1").Value = "RTDO V.1"

# Swap / update player names
$ws.Range("H2").Value = "Gonzo"

$ws.Range("G3").Value = "Coquina"
$ws.Range("H3").Value = "Puche"

$ws.Range("G4").Value = "Ruso"
$ws.Range("H4").Value = "Lope"
